$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 886 (shifts rows 886:977 down to 887:978,
# and Excel copies the formatting of the row above into the new row, including
# the date-style cell in column D).
$ws.Rows.Item(886).Insert()

# Populate the newly inserted row 886 with the new data point.
$ws.Cells.Item(886, 1).Value = 3
$ws.Cells.Item(886, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(886, 3).Value = "Coquimbo"
$ws.Cells.Item(886, 4).Value2 = 45166
$ws.Cells.Item(886, 5).Value = 5
$ws.Cells.Item(886, 6).Value = 100112045
$ws.Cells.Item(886, 7).Value = "Zapallo"
$ws.Cells.Item(886, 8).Value = "Camote"
$ws.Cells.Item(886, 9).Value = "1a (guarda)"
$ws.Cells.Item(886, 10).Value = 120
$ws.Cells.Item(886, 11).Value = 1100
$ws.Cells.Item(886, 12).Value = 1100
$ws.Cells.Item(886, 13).Value = 1100
$ws.Cells.Item(886, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(886, 15).Value = "Provincia de Talca"
$ws.Cells.Item(886, 16).Value = 1100
$ws.Cells.Item(886, 17).Value = 1
$ws.Cells.Item(886, 18).Value = "Hortaliza"
